# Daily attendance processing - 2025-10-24 22:20:24
#
# Re-orders the comma-separated contributor list in column G ("Recorded By")
# for the session-analysis rows whose recorder list currently lists the
# human contributor(s) before "System"/"system" - flips it so "System"
# comes first (or is pushed toward the tail, for the 3-name case), matching
# the day's reprocessed attendance log ordering.
#
# Exact value substitutions applied (order-preserving list reversal):
#   "backup@backdoor.com, System, system" -> "system, System, backup@backdoor.com"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
#
# Rows outside this known set (e.g. "System, admin@admin.com",
# "dnasr281@gmail.com, admin@admin.com", single-name entries, blanks) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

# Rows with a "Recorded By" value affected by this reprocessing pass.
$targetRows = @(
    2,3,4,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,
    29,30,31,32,33,35,37,38,39,40,41,42,44,45,46,47,48,49,
    56,57,58,59,60,62,64,65,66,67,68,69,71,72,73,74,75,76,
    83,84,85,86,87,88,89,93,95,96,97,99,102,
    109,110,111,112,113,114,115,119,121,122,123,125,128,
    135,136,137,138,139,140,141,145,147,148,149,151,154
)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value()
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
